$d = $word.ActiveDocument

# 1. Remove the "Map - A map of the community, could be an embedded Google
#    map" bullet entirely (whole paragraph, including its paragraph mark).
$mapRange = $d.Range(0, 0)
[void]$mapRange.Find.Execute("Map", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$mapPara = $mapRange.Paragraphs(1)
$mapPara.Range.Delete()

# 2. "Important locations and links if available" -> "Important locations
#    and links " (drop the trailing "if available").
[void]$d.Content.Find.Execute("Important locations and links if available", $false, $false, $false, $false, $false, $true, 1, $false, "Important locations and links ", 2)

# 3. "Significant people and links if available" -> split the run after
#    "links " and wrap "if available" with the (relocated) _GoBack bookmark,
#    matching how Word leaves its "last edit" marker behind after typing.
$prefix = $d.Range(0, 0)
[void]$prefix.Find.Execute("Significant people and links ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$fullPara = $d.Range(0, 0)
[void]$fullPara.Find.Execute("Significant people and links if available", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$ifAvailable = $d.Range($prefix.End, $fullPara.End)
$d.Bookmarks.Add("_GoBack", $ifAvailable) | Out-Null

Write-Host "Edits applied."
